$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 9501.5
$ws.Range("I76").Value = 9501.5
$ws.Range("K76").Value = 9501.5
$ws.Range("M76").Value = -9186.5
$ws.Range("H79").Value = 9501.5
$ws.Range("I79").Value = 9501.5
$ws.Range("K79").Value = 9501.5
$ws.Range("M79").Value = -8409.5
$ws.Range("H80").Value = 1974.75
$ws.Range("J80").Value = 2333
$ws.Range("L80").Value = 6999
$ws.Range("N80").Value = -8995
$ws.Range("H83").Value = 1974.75
$ws.Range("J83").Value = 2333
$ws.Range("L83").Value = 20997
$ws.Range("N83").Value = -30981
$ws.Range("H86").Value = 2999
$ws.Range("I86").Value = 2999
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2999
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -1876
$ws.Range("H89").Value = 2999
$ws.Range("I89").Value = 2999
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 14995
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -9379
$ws.Range("H92").Value = 721.82355
$ws.Range("I92").Value = 721.82355
$ws.Range("K92").Value = 721.82355
$ws.Range("M92").Value = 526.17645
$ws.Range("H107").Value = 1104.5
$ws.Range("I107").Value = 1087.4
$ws.Range("K107").Value = 1087.4
$ws.Range("M107").Value = 832.5999999999999
$ws.Range("H137").Value = 2065
$ws.Range("J137").Value = 1859.8125
$ws.Range("L137").Value = 5579.4375
$ws.Range("N137").Value = -10679.4375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 498.63635
$ws.Range("I4").Value = 287.22223
$ws.Range("K4").Value = 287.22223
$ws.Range("M4").Value = -171.22223
$ws.Range("H63").Value = 3731.7
$ws.Range("J63").Value = 3590.7778
$ws.Range("L63").Value = 3590.7778
$ws.Range("N63").Value = -4962.7778
$ws.Range("H66").Value = 3731.7
$ws.Range("J66").Value = 3590.7778
$ws.Range("L66").Value = 17953.889
$ws.Range("N66").Value = -24817.889
$ws.Range("H122").Value = 1767.2
$ws.Range("I122").Value = 1385.5714
$ws.Range("J122").Value = 2657.6667
$ws.Range("K122").Value = 4156.7142
$ws.Range("L122").Value = 7973.000100000001
$ws.Range("M122").Value = -1706.7142
$ws.Range("N122").Value = -12873.0001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("H76").Value = 22500
$ws.Range("J76").Value = 22500
$ws.Range("L76").Value = 22500
$ws.Range("N76").Value = -23130
$ws.Range("H79").Value = 22500
$ws.Range("J79").Value = 22500
$ws.Range("L79").Value = 22500
$ws.Range("N79").Value = -24684
$ws.Range("H105").Value = 2450
$ws.Range("I105").Value = 2713.4443
$ws.Range("J105").Value = 1659.6666
$ws.Range("K105").Value = 2713.4443
$ws.Range("L105").Value = 1659.6666
$ws.Range("M105").Value = -966.4443000000001
$ws.Range("N105").Value = -5153.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 73913.5
$ws.Range("J31").Value = 6884.5
$ws.Range("L31").Value = 6884.5
$ws.Range("N31").Value = -7474.5
$ws.Range("H34").Value = 73913.5
$ws.Range("J34").Value = 6884.5
$ws.Range("L34").Value = 6884.5
$ws.Range("N34").Value = -7288.5
$ws.Range("H99").Value = 2505.1177
$ws.Range("I99").Value = 2422.1538
$ws.Range("J99").Value = 2774.75
$ws.Range("K99").Value = 2422.1538
$ws.Range("L99").Value = 2774.75
$ws.Range("M99").Value = -924.1538
$ws.Range("N99").Value = -5770.75
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 1594.2222
$ws.Range("I105").Value = 1606
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1606
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 141
$ws.Range("N105").Value = -4994
$ws.Range("H126").Value = 2505.1177
$ws.Range("I126").Value = 2422.1538
$ws.Range("J126").Value = 2774.75
$ws.Range("K126").Value = 7266.4614
$ws.Range("L126").Value = 8324.25
$ws.Range("M126").Value = -4796.4614
$ws.Range("N126").Value = -13264.25
$ws.Range("H132").Value = 2083.2126
$ws.Range("I132").Value = 1850.119
$ws.Range("J132").Value = 4041.2
$ws.Range("K132").Value = 5550.357
$ws.Range("L132").Value = 12123.6
$ws.Range("M132").Value = -3020.357
$ws.Range("N132").Value = -17183.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 53072.57
$ws.Range("I4").Value = 849.5625
$ws.Range("J4").Value = 220186.2
$ws.Range("K4").Value = 2548.6875
$ws.Range("L4").Value = 660558.6000000001
$ws.Range("M4").Value = -2436.6875
$ws.Range("N4").Value = -660782.6000000001
$ws.Range("H114").Value = 2430
$ws.Range("J114").Value = 3999.3333
$ws.Range("L114").Value = 11997.9999
$ws.Range("N114").Value = -18505.9999
$ws.Range("H139").Value = 2381.8
$ws.Range("I139").Value = 2272.7058
$ws.Range("K139").Value = 6818.117400000001
$ws.Range("M139").Value = -1678.117400000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 37072036
$ws.Range("I18").Value = 55558056
$ws.Range("K18").Value = 55558056
$ws.Range("M18").Value = -55557763
$ws.Range("H122").Value = 3092.75
$ws.Range("I122").Value = 2452.5386
$ws.Range("K122").Value = 7357.6158
$ws.Range("M122").Value = -4907.6158

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 16720000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H82").Value = 1723.6
$ws.Range("I82").Value = 2075.5
$ws.Range("J82").Value = 1489
$ws.Range("K82").Value = 2075.5
$ws.Range("L82").Value = 1489
$ws.Range("M82").Value = -1714.5
$ws.Range("N82").Value = -2211
$ws.Range("H85").Value = 1723.6
$ws.Range("I85").Value = 2075.5
$ws.Range("J85").Value = 1489
$ws.Range("K85").Value = 2075.5
$ws.Range("L85").Value = 1489
$ws.Range("M85").Value = -827.5
$ws.Range("N85").Value = -3985
$ws.Range("H122").Value = 5391.778
$ws.Range("I122").Value = 4850.615
$ws.Range("J122").Value = 6798.8
$ws.Range("K122").Value = 14551.845
$ws.Range("L122").Value = 20396.4
$ws.Range("M122").Value = -12101.845
$ws.Range("N122").Value = -25296.4
$ws.Range("H132").Value = 4904.857
$ws.Range("I132").Value = 4151.727
$ws.Range("K132").Value = 12455.181
$ws.Range("M132").Value = -9925.181

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 50000
$ws.Range("J27").Value = 50000
$ws.Range("L27").Value = 50000
$ws.Range("N27").Value = -50138
$ws.Range("H81").Value = 6638.724
$ws.Range("I81").Value = 8855.77
$ws.Range("J81").Value = 4837.375
$ws.Range("K81").Value = 17711.54
$ws.Range("L81").Value = 9674.75
$ws.Range("M81").Value = -16650.54
$ws.Range("N81").Value = -11796.75
$ws.Range("H84").Value = 6638.724
$ws.Range("I84").Value = 8855.77
$ws.Range("J84").Value = 4837.375
$ws.Range("K84").Value = 88557.70000000001
$ws.Range("L84").Value = 48373.75
$ws.Range("M84").Value = -83253.70000000001
$ws.Range("N84").Value = -58981.75
$ws.Range("H122").Value = 2184.0557
$ws.Range("I122").Value = 2127.8667
$ws.Range("K122").Value = 6383.6001
$ws.Range("M122").Value = -3933.6001
$ws.Range("H126").Value = 1981.9474
$ws.Range("I126").Value = 1866.2307
$ws.Range("K126").Value = 5598.6921
$ws.Range("M126").Value = -3128.6921
$ws.Range("H132").Value = 2280.7778
$ws.Range("I132").Value = 2026.2162
$ws.Range("J132").Value = 3458.125
$ws.Range("K132").Value = 6078.6486
$ws.Range("L132").Value = 10374.375
$ws.Range("M132").Value = -3548.6486
$ws.Range("N132").Value = -15434.375
